$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/parent-organization-hierarchy-level-description"
$wsMeta.Range("B3").Value = "8.0.0"
$wsMeta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMeta.Range("B9").Value = "LinuxForHealth Team"

$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("AI2").Value = ""
$wsElem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/parent-organization-hierarchy-level-description"
